$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update file names to include "_test" suffix
$ws.Range("A2").Value = "Ctrl_plate1_leukocytes_test.fcs"
$ws.Range("A3").Value = "P15_D1_leukocytes_test.fcs"

# Update the selected cell to match the final state
$ws.Range("D9").Select()
